$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-25: columns B,C,D,E,F,I,J,K,L,M,N
$data = @(
    @(1.02, 1.034063234460032, 1.060760384964664, 1.03329305325899, 1.064979749507142, 1.048272377781091, 1.039184532615147, 1.063486287743405, 1.036095820211693, 1.067694225211488, 1.040660293847593),
    @(1.02, 1.034964414879229, 1.061479701356084, 1.034057817870815, 1.065832694152819, 1.048541584049146, 1.039728932610811, 1.064020332827844, 1.0366696529346, 1.068362385759404, 1.041205466953693),
    @(1.02, 1.035548081585984, 1.061945348833855, 1.034553512338095, 1.066385210360045, 1.04871467736394, 1.040081104745508, 1.06436542383891, 1.037041140567214, 1.06879468697315, 1.041558139213217),
    @(1.02, 1.035793583768162, 1.062141153632356, 1.034762102010859, 1.066617631031238, 1.048787181408671, 1.040229135040257, 1.064510385823956, 1.037197355777031, 1.068976414558635, 1.041706379727964),
    @(1.02, 1.035834812176139, 1.06217403280325, 1.034797136801907, 1.06665666381134, 1.048799339630828, 1.040253988605892, 1.064534718839831, 1.037223587416093, 1.069006926682577, 1.041731268588513),
    @(1.02, 1.03555136149399, 1.061947965006169, 1.034556298740674, 1.066388315417278, 1.048715647206345, 1.040083082824757, 1.064367361278719, 1.037043227759668, 1.068797115274427, 1.041560120101566),
    @(1.02, 1.034367679489501, 1.061003438538364, 1.033551334064318, 1.06526788003944, 1.048363584835726, 1.039368533355099, 1.063666867696485, 1.03628971180437, 1.067920041287431, 1.0408445558897),
    @(1.02, 1.032286091443423, 1.059340681061949, 1.03178696361409, 1.06329823232911, 1.047734806522946, 1.038108755683676, 1.062428952378195, 1.034963346387475, 1.066374257820361, 1.039582989189485),
    @(1.02, 1.030901269640955, 1.058233364299404, 1.030615181462348, 1.061988398683386, 1.047310019222725, 1.037268530009721, 1.06160136217535, 1.034080130330965, 1.065343643049064, 1.038741570298698),
    @(1.02, 1.030302330187845, 1.057754186697732, 1.030108864898506, 1.061422022238735, 1.047124762464797, 1.036904627537123, 1.061242474041194, 1.033697946151699, 1.064897371780188, 1.038377151042849),
    @(1.02, 1.030079963374993, 1.057576245228481, 1.029920959095743, 1.061211765252056, 1.047055752119926, 1.036769447098055, 1.061109087751019, 1.033556025267885, 1.064731606968177, 1.038241778632052),
    @(1.02, 1.030127656997055, 1.057614412149832, 1.029961258140634, 1.061256860611389, 1.047070564009302, 1.036798444235714, 1.061137703128957, 1.033586465980733, 1.064767164050384, 1.038270816948972),
    @(1.02, 1.030283947109218, 1.057739477043579, 1.030093329218091, 1.061404639863632, 1.047119062077657, 1.036893453695325, 1.061231449898225, 1.033686214118075, 1.064883669596497, 1.038365961332915),
    @(1.02, 1.030380256682462, 1.057816539854485, 1.030174724174603, 1.061495707549818, 1.047148917158288, 1.036951990751517, 1.061289199912908, 1.033747677484436, 1.064955452579498, 1.03842458151844),
    @(1.02, 1.030941034075939, 1.058265172157502, 1.030648806759364, 1.062026004018223, 1.047322286333404, 1.037292679438425, 1.061625169213424, 1.034105500089816, 1.065373260552932, 1.038765754022362),
    @(1.02, 1.031292982493077, 1.058546668279256, 1.030946474550333, 1.062358857753644, 1.047430683072737, 1.037506363979478, 1.061835771255214, 1.034330021622678, 1.065635339265521, 1.038979742019973),
    @(1.02, 1.031498335272557, 1.058710888831711, 1.031120202547608, 1.062553081829859, 1.047493781544815, 1.037630994844319, 1.061958560049413, 1.034461005684074, 1.065788204449184, 1.039104549874937),
    @(1.02, 1.031568366630025, 1.058766888600888, 1.031179456798857, 1.062619320079953, 1.047515274842532, 1.03767348940857, 1.062000419027419, 1.034505671976336, 1.065840327330431, 1.039147104786343),
    @(1.02, 1.031255214797124, 1.058516463436286, 1.030914526906125, 1.062323137811594, 1.047419066312305, 1.037483438440735, 1.061813180999271, 1.034305930049802, 1.065607220782078, 1.038956783924335),
    @(1.02, 1.030237920631944, 1.057702647254449, 1.030054433073606, 1.061361119218453, 1.047104786057874, 1.036865476068873, 1.06120384598881, 1.03365683968268, 1.064849361599616, 1.038337943975023),
    @(1.02, 1.029598921094876, 1.057191238164081, 1.029514600107283, 1.060756957646926, 1.046906041912309, 1.03647687629477, 1.060820274960293, 1.033248958966609, 1.064372867370574, 1.037948792344668),
    @(1.02, 1.029937608277183, 1.057462319699483, 1.029800685959614, 1.061077168415928, 1.047011508104957, 1.036682885944862, 1.061023656195802, 1.033465162241806, 1.064625465310873, 1.038155094552093),
    @(1.02, 1.031272280185633, 1.058530111616414, 1.030928962352029, 1.062339277880847, 1.047424315820371, 1.037493797528814, 1.061823388723783, 1.034316815919154, 1.065619926317327, 1.038967157723506),
    @(1.02, 1.032823725392459, 1.059770341594429, 1.0322423153541, 1.063806865571129, 1.047898351476906, 1.038434508945279, 1.062749396581427, 1.035306067565845, 1.066773902570599, 1.039909205058084)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("I$row").Value = $vals[5]
    $ws.Range("J$row").Value = $vals[6]
    $ws.Range("K$row").Value = $vals[7]
    $ws.Range("L$row").Value = $vals[8]
    $ws.Range("M$row").Value = $vals[9]
    $ws.Range("N$row").Value = $vals[10]
}
